$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing blank row (row 6); the used range shrinks to A1:V5
$ws.Rows.Item(6).Delete()

# Correct two shared-string lookups that pointed at stale values
# NivelAtual for row 2 (M2): "IX" -> "I"
$ws.Cells.Item(2, 13).Value = "I"
# GrauNovo for row 4 (V4): "K" -> "G"
$ws.Cells.Item(4, 22).Value = "G"

# Populate row 5, which previously held only blank, styled cells, with a
# duplicate of row 2's submission record
$ws.Cells.Item(5, 1).Value = "qJPar8"
$ws.Cells.Item(5, 2).Value = "jBVv5Q"
$ws.Cells.Item(5, 3).Value = 45425.798657407409
$ws.Cells.Item(5, 4).Value = "andre.amorim@planejamento.mg.gov.br"
$ws.Cells.Item(5, 5).Value = "teste"
$ws.Cells.Item(5, 6).Value = "André"
$ws.Cells.Item(5, 7).Value = 1
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 44329
$ws.Cells.Item(5, 11).Value = 44692
$ws.Cells.Item(5, 12).Value = "EPPGG"
$ws.Cells.Item(5, 13).Value = "I"
$ws.Cells.Item(5, 14).Value = "EPPGG"
$ws.Cells.Item(5, 15).Value = "A"
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 0
$ws.Cells.Item(5, 18).Value = 45426
$ws.Cells.Item(5, 19).Value = 45790
$ws.Cells.Item(5, 20).Value = 45413
$ws.Cells.Item(5, 21).Value = "I"
$ws.Cells.Item(5, 22).Value = "B"

# Rows 1-4 lose their explicit 12.75pt height and fall back to the sheet
# default (15.75pt, no customHeight flag)
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
